$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Update the "sudan" sheet dictionary data
# ---------------------------------------------------------------------------
$sudan = $wb.Worksheets.Item("sudan")

# Clear out the old data range (it had 24 rows, A1:C24) before writing new data
$sudan.Range("A1:C40").ClearContents()

$sudanData = @(
    @("indicator", "variable", "label"),
    @(1.1, "I1_sec_day", "Feeling safe at day"),
    @(1.1, "I1_sec_night", "Feeling safe at night"),
    @(1.1, "I1_sec_inc", "Experience security incident"),
    @(1.1, "I1_sec_rep", "Report security incident"),
    @(1.1, "I1_SDG_16.1.4", "Feeling safe at day and night"),
    @(2.1, "I3_no_borrow", "Borrowing for food"),
    @(2.1, "I3_pay_food", "Ability to pay for food"),
    @(2.1, "I3_DS_2.1.2", "Food insecurity scale"),
    @(2.2, "I4_hous_ownership", "Own house"),
    @(2.2, "I4_hous_water", "Improved water "),
    @(2.2, "I4_hous_toilet", "Improved sanitation"),
    @(2.2, "I4_hous_overcrowd", "Overcrowded housing"),
    @(2.2, "I4_hous_permanent", "Permanent housing structures"),
    @(2.2, "I4_SDG_11.1.1", "Living outside of slums"),
    @(2.3, "I5_DS_2.1.8", "Satisfied with health facilities"),
    @(2.3, "I5_med_dist", "Distance to health facility"),
    @(2.4, "I6_edu_dist", "Distance to school"),
    @(2.4, "I6_SDG_4.1.2", "Ever in school"),
    @(2.4, "I6_educ_child", "Child in school"),
    @(3.1, "I7_job_unemploy", "Employed"),
    @(3.1, "I7_SDG_8.5.2", "Unemployment"),
    @(3.2, "I8_econ_account", "Bank account"),
    @(3.2, "I8_econ_market", "Access to market"),
    @(3.2, "I8_SDG_1.2.1", "Below 1.9 USD Poverty Line"),
    @(3.2, "I8_poor32", "Below 3.2 USD Poverty Line"),
    @(4.1, "I9_SDG_1.4.2", "Security of tenure"),
    @(4.1, "I9_hlp_access", "Access to compensation"),
    @(4.1, "I9_hlp_doc", "Documentation"),
    @(4.1, "I9_hlp_own", "Ownership over property"),
    @(5.1, "I10_doc_birth", "Birth certificate"),
    @(5.1, "I10_DS_5.1.1", "Possession of ID")
)

for ($i = 0; $i -lt $sudanData.Length; $i++) {
    $r = $i + 1
    $row = $sudanData[$i]
    $sudan.Cells.Item($r, 1).Value = $row[0]
    $sudan.Cells.Item($r, 2).Value = $row[1]
    $sudan.Cells.Item($r, 3).Value = $row[2]
}

# ---------------------------------------------------------------------------
# 2. Update view/selection state for each sheet
# ---------------------------------------------------------------------------
$hargeisa = $wb.Worksheets.Item("hargeisa")
$nigeria = $wb.Worksheets.Item("nigeria")

$hargeisa.Activate()
$hargeisa.Range("A29:XFD29").Select()

$nigeria.Activate()
$nigeria.Range("A23:XFD23").Select()

$sudan.Activate()
$sudan.Range("D27").Select()
